$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (Total) summary sheet: insert a new row for 2022-Q4
#    at the top of the data (row 2), pushing the existing quarters down by
#    one row, and keep the running index in column A sequential.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Copy the column-A style (s="2") from the row below onto the new row.
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122) # xlPasteFormats

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 15
$total.Cells.Item(2,4).Value = 2.09

for ($r = 3; $r -le 10; $r++) {
    $total.Cells.Item($r,1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" worksheet right after "总计", holding the
#    fund-holding breakdown for the new quarter.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row (bold, centered, top-aligned, thin border) — matches the style
# used on every other quarter sheet's header row.
$header = $q4.Range("A1:H1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108 # xlCenter
$header.VerticalAlignment = -4160   # xlTop
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

$rows = @(
    @(0, "009010", "华夏兴阳一年持有期混合", "27.18", "90.65", "2.32", "0.6306", 9),
    @(1, "009846", "富兰克林国海港股通远见价值混合", "15.17", "90.16", "2.98", "0.4521", 8),
    @(2, "005888", "华夏新兴消费混合A", "7.99", "88.30", "2.98", "0.2381", 7),
    @(3, "202801", "南方全球精选配置（QDII-FOF）", "17.02", "32.64", "1.18", "0.2008", 8),
    @(4, "005889", "华夏新兴消费混合C", "4.95", "88.30", "2.98", "0.1475", 7),
    @(5, "007182", "万家沪港深蓝筹混合A", "3.44", "92.13", "4.09", "0.1407", 3),
    @(6, "013009", "万家港股通精选混合A", "2.56", "84.90", "3.76", "0.0963", 3),
    @(7, "007183", "万家沪港深蓝筹混合C", "0.86", "92.13", "4.09", "0.0352", 3),
    @(8, "005646", "中海沪港深多策略灵活配置混合", "0.79", "105.93", "4.35", "0.0344", 9),
    @(9, "013010", "万家港股通精选混合C", "0.77", "84.90", "3.76", "0.0290", 3),
    @(10, "005143", "中融沪港深大消费主题灵活配置混合C", "0.69", "77.26", "4.19", "0.0289", 7),
    @(11, "013767", "平安价值回报混合A", "0.83", "91.57", "2.82", "0.0234", 7),
    @(12, "009140", "永赢竞争力精选混合", "0.60", "94.52", "3.62", "0.0217", 6),
    @(13, "005142", "中融沪港深大消费主题灵活配置混合A", "0.31", "77.26", "4.19", "0.0130", 7),
    @(14, "013768", "平安价值回报混合C", "0.04", "91.57", "2.82", "0.0011", 7)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r,1).Value = $row[0]
    $q4.Cells.Item($r,1).Font.Bold = $true
    $q4.Cells.Item($r,1).HorizontalAlignment = -4108
    $q4.Cells.Item($r,1).VerticalAlignment = -4160
    $q4.Cells.Item($r,1).Borders.LineStyle = 1
    $q4.Cells.Item($r,1).Borders.Weight = 2

    $q4.Cells.Item($r,2).NumberFormat = "@"
    $q4.Cells.Item($r,2).Value = $row[1]
    $q4.Cells.Item($r,3).Value = $row[2]
    $q4.Cells.Item($r,4).NumberFormat = "@"
    $q4.Cells.Item($r,4).Value = $row[3]
    $q4.Cells.Item($r,5).NumberFormat = "@"
    $q4.Cells.Item($r,5).Value = $row[4]
    $q4.Cells.Item($r,6).NumberFormat = "@"
    $q4.Cells.Item($r,6).Value = $row[5]
    $q4.Cells.Item($r,7).NumberFormat = "@"
    $q4.Cells.Item($r,7).Value = $row[6]
    $q4.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3. Restore the originally-selected tab (2020-Q4, now the last sheet).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
